# Weekly update: a new data row is inserted right before the current row 86
# (shifting rows 86..163 down to 87..164) and populated with a new weekly
# observation (same commodity/variety/quality/origin as the old row 86, but
# a new date and a new Volumen figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86 — everything currently at row 86 downward shifts
# down by one (old row 86 -> 87, ..., old row 163 -> 164).
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly entry.
$ws.Range("A86").Value = 8
$ws.Range("B86").Value = "Terminal La Palmera de La Serena"
$ws.Range("C86").Value = "Coquimbo"
$ws.Range("D86").Value = 44566
$ws.Range("E86").Value = 4
$ws.Range("F86").Value = 100112031
$ws.Range("G86").Value = "Poroto verde"
$ws.Range("H86").Value = "Magnum"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 25000
$ws.Range("L86").Value = 26000
$ws.Range("M86").Value = 25500
$ws.Range("N86").Value = '$/malla 25 kilos'
$ws.Range("O86").Value = 'Provincia de Limarí'
$ws.Range("P86").Value = 1020
$ws.Range("Q86").Value = 25
$ws.Range("R86").Value = "Hortaliza"
